# Progress in sheet3 parse
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("End Systems")

# Fill in the running index column D (rows 2-13) with 1..12
for ($i = 2; $i -le 13; $i++) {
    $ws.Cells.Item($i, 4).Value = $i - 1
}

# Add the new sample/legend row 17: column headers 0..13 across A..N
for ($col = 1; $col -le 14; $col++) {
    $ws.Cells.Item(17, $col).Value = $col - 1
}

# Move the active selection to A2
$ws.Range("A2").Select()
